# Auto-generated Excel COM-interop script applying the market-price refresh diff.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H, I, J, K, L, M, N)
# for the rows whose backing market data changed in this scheduled refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2474.75
$ws.Range("J17").Value = 2474.75
$ws.Range("L17").Value = 7424.25
$ws.Range("N17").Value = -7760.25

$ws.Range("H62").Value = 14282.381
$ws.Range("I62").Value = 11114.728
$ws.Range("J62").Value = 17766.8
$ws.Range("K62").Value = 11114.728
$ws.Range("L62").Value = 17766.8
$ws.Range("M62").Value = -10490.728
$ws.Range("N62").Value = -19014.8

$ws.Range("H65").Value = 14282.381
$ws.Range("I65").Value = 11114.728
$ws.Range("J65").Value = 17766.8
$ws.Range("K65").Value = 55573.64
$ws.Range("L65").Value = 88834
$ws.Range("M65").Value = -52453.64
$ws.Range("N65").Value = -95074

$ws.Range("H100").Value = 2777.7778
$ws.Range("I100").Value = 1889.5
$ws.Range("J100").Value = 3888.125
$ws.Range("K100").Value = 1889.5
$ws.Range("L100").Value = 3888.125
$ws.Range("M100").Value = -1348.5
$ws.Range("N100").Value = -4970.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15765.634
$ws.Range("I32").Value = 15258.363
$ws.Range("J32").Value = 23577.6
$ws.Range("K32").Value = 15258.363
$ws.Range("L32").Value = 23577.6
$ws.Range("M32").Value = -14971.363
$ws.Range("N32").Value = -24151.6

$ws.Range("H61").Value = 2927.889
$ws.Range("I61").Value = 2680.6
$ws.Range("J61").Value = 3237
$ws.Range("K61").Value = 2680.6
$ws.Range("L61").Value = 3237
$ws.Range("M61").Value = -2468.6
$ws.Range("N61").Value = -3661

$ws.Range("H74").Value = 36335.832
$ws.Range("I74").Value = 39964.777
$ws.Range("J74").Value = 3675.3333
$ws.Range("K74").Value = 39964.777
$ws.Range("L74").Value = 3675.3333
$ws.Range("M74").Value = -39090.777
$ws.Range("N74").Value = -5423.3333

$ws.Range("H77").Value = 36335.832
$ws.Range("I77").Value = 39964.777
$ws.Range("J77").Value = 3675.3333
$ws.Range("K77").Value = 199823.885
$ws.Range("L77").Value = 18376.6665
$ws.Range("M77").Value = -195455.885
$ws.Range("N77").Value = -27112.6665

$ws.Range("H122").Value = 2166.923
$ws.Range("I122").Value = 1972.5
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 5917.5
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -3467.5
$ws.Range("N122").Value = -18400

$ws.Range("H132").Value = 27507.924
$ws.Range("I132").Value = 29466.916
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 88400.74800000001
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -85870.74800000001
$ws.Range("N132").Value = -17060

$ws.Range("H136").Value = 2927.889
$ws.Range("I136").Value = 2680.6
$ws.Range("J136").Value = 3237
$ws.Range("K136").Value = 8041.799999999999
$ws.Range("L136").Value = 9711
$ws.Range("M136").Value = -5491.799999999999
$ws.Range("N136").Value = -14811

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 40011.668
$ws.Range("J38").Value = 40011.668
$ws.Range("L38").Value = 40011.668
$ws.Range("N38").Value = -40843.668

$ws.Range("H80").Value = 451
$ws.Range("I80").Value = 508
$ws.Range("J80").Value = 299
$ws.Range("K80").Value = 508
$ws.Range("L80").Value = 299
$ws.Range("M80").Value = 490
$ws.Range("N80").Value = -2295

$ws.Range("H83").Value = 451
$ws.Range("I83").Value = 508
$ws.Range("J83").Value = 299
$ws.Range("K83").Value = 2540
$ws.Range("L83").Value = 1495
$ws.Range("M83").Value = 2452
$ws.Range("N83").Value = -11479

$ws.Range("H105").Value = 2435.074
$ws.Range("I105").Value = 2365.9546
$ws.Range("J105").Value = 2739.2
$ws.Range("K105").Value = 2365.9546
$ws.Range("L105").Value = 2739.2
$ws.Range("M105").Value = -618.9546
$ws.Range("N105").Value = -6233.2

$ws.Range("H134").Value = 3267.5293
$ws.Range("I134").Value = 2710.6428
$ws.Range("K134").Value = 8131.928400000001
$ws.Range("M134").Value = -5596.928400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3752.6365
$ws.Range("I31").Value = 3277.4
$ws.Range("J31").Value = 4771
$ws.Range("K31").Value = 3277.4
$ws.Range("L31").Value = 4771
$ws.Range("M31").Value = -2982.4
$ws.Range("N31").Value = -5361

$ws.Range("H34").Value = 3752.6365
$ws.Range("I34").Value = 3277.4
$ws.Range("J34").Value = 4771
$ws.Range("K34").Value = 3277.4
$ws.Range("L34").Value = 4771
$ws.Range("M34").Value = -3075.4
$ws.Range("N34").Value = -5175

$ws.Range("H58").Value = 48178.5
$ws.Range("I58").Value = 58245.945
$ws.Range("J58").Value = 2875
$ws.Range("K58").Value = 58245.945
$ws.Range("L58").Value = 2875
$ws.Range("M58").Value = -58042.945
$ws.Range("N58").Value = -3281

$ws.Range("H132").Value = 3248.3333
$ws.Range("I132").Value = 3148.5217
$ws.Range("J132").Value = 3576.2856
$ws.Range("K132").Value = 9445.5651
$ws.Range("L132").Value = 10728.8568
$ws.Range("M132").Value = -6915.5651
$ws.Range("N132").Value = -15788.8568

$ws.Range("H136").Value = 48178.5
$ws.Range("I136").Value = 58245.945
$ws.Range("J136").Value = 2875
$ws.Range("K136").Value = 174737.835
$ws.Range("L136").Value = 8625
$ws.Range("M136").Value = -172187.835
$ws.Range("N136").Value = -13725

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 412.83334
$ws.Range("I97").Value = 282
$ws.Range("J97").Value = 456.44446
$ws.Range("K97").Value = 846
$ws.Range("L97").Value = 1369.33338
$ws.Range("M97").Value = -350
$ws.Range("N97").Value = -2361.33338

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 36054.965
$ws.Range("I132").Value = 46065.39
$ws.Range("J132").Value = 3163.5715
$ws.Range("K132").Value = 138196.17
$ws.Range("L132").Value = 9490.7145
$ws.Range("M132").Value = -135666.17
$ws.Range("N132").Value = -14550.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4017.6365
$ws.Range("I68").Value = 3824.25
$ws.Range("K68").Value = 3824.25
$ws.Range("M68").Value = -3075.25

$ws.Range("H71").Value = 4017.6365
$ws.Range("I71").Value = 3824.25
$ws.Range("K71").Value = 19121.25
$ws.Range("M71").Value = -15377.25

$ws.Range("H100").Value = 4169
$ws.Range("I100").Value = 3804.818
$ws.Range("J100").Value = 4569.6
$ws.Range("K100").Value = 3804.818
$ws.Range("L100").Value = 4569.6
$ws.Range("M100").Value = -3263.818
$ws.Range("N100").Value = -5651.6

$ws.Range("H132").Value = 42862
$ws.Range("I132").Value = 51481.4
$ws.Range("J132").Value = 6947.8335
$ws.Range("K132").Value = 154444.2
$ws.Range("L132").Value = 20843.5005
$ws.Range("M132").Value = -151914.2
$ws.Range("N132").Value = -25903.5005

Write-Output "Sheets updated via scheduled runner"

